# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds strikeout counts that were previously derived from a
# "Strike#" source and must be regenerated from the corrected K values.
# Apply the corrected values row by row (data starts at row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 2
    25 = 3
    26 = 3
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    31 = 3
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 2
    37 = 2
    38 = 0
    39 = 0
    40 = 1
    41 = 2
    42 = 0
    43 = 0
    44 = 4
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 2
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    56 = 0
    57 = 0
    58 = 3
    59 = 1
    60 = 0
    61 = 3
    62 = 2
    63 = 0
    64 = 2
    65 = 3
    66 = 2
    67 = 1
    68 = 1
    69 = 1
    70 = 1
    71 = 0
    72 = 2
    73 = 2
    74 = 0
    75 = 1
    76 = 2
    77 = 0
    78 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
